$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New dictionary entry: "Neonatal Deaths" appended as the next row (58) of
# the disease/threshold table, matching the formatting (fill + border +
# wrap) of the preceding rows.
$ws.Range("A57:E57").Copy()
$ws.Range("A58:E58").PasteSpecial(-4122)

$ws.Range("A58").Value = "Neonatal Deaths"
$ws.Range("B58").Value = "1 death"
$ws.Range("C58").Value = "Verify, ensure a case review is performed and report as soon as possible/ Vérifier,  s'assurez qu'une revue des cas a été réalisée et notifier le plus tôt possible"
$ws.Range("E58").Value = "Death"
# D58 (Comments) is left blank, matching the rest of the table.

# "Neonatal Deaths" is long enough to need wrapping in column A (unlike the
# preceding rows), and the row grows to the same wrapped height as the rest
# of the table.
$ws.Range("A58").WrapText = $true
$ws.Rows.Item(58).RowHeight = 43.5

# Scroll the view down to the newly-added row and select the cell below it,
# reflecting where the author ended up after the edit.
$ws.Range("A59").Select()
$excel.ActiveWindow.ScrollRow = 55
